# Insert one new weekly price record for "Orégano" into the dataset.
# The new record belongs right after the current header block of rows and
# before the existing row 49, so every following row shifts down by one
# (old row 49 -> new row 50, ..., old row 158 -> new row 159).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 49..158 down by inserting a new blank row at position 49.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new observation.
$ws.Cells.Item(49, 1).Value  = 6
$ws.Cells.Item(49, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(49, 3).Value  = "Metropolitana"
$ws.Cells.Item(49, 4).Value  = 44614
$ws.Cells.Item(49, 5).Value  = 13
$ws.Cells.Item(49, 6).Value  = 100112029
$ws.Cells.Item(49, 7).Value  = "Orégano"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 34
$ws.Cells.Item(49, 11).Value = 9000
$ws.Cells.Item(49, 12).Value = 10000
$ws.Cells.Item(49, 13).Value = 9441
$ws.Cells.Item(49, 14).Value = "`$/docena de atados"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 3147
$ws.Cells.Item(49, 17).Value = 3
$ws.Cells.Item(49, 18).Value = "Hortaliza"
